# Edit script for 3busCDF.xlsx
# Implements:
#  - Bus 1 (row 3): H/I load values 1500/750 -> 0/0 (style becomes General, like col J)
#  - Bus 2 (row 4): H/I load values 900/500 -> 1500/750 (style unchanged)
#  - Insert a new "Bus 3" row (new row 5) with the load values 900/500 that used
#    to belong to Bus 2, shifting all subsequent rows down by one
#  - Update the active-cell selection to I7 (was I8) to reflect the shifted layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Bus 1): zero-out the load, switching number format to General ---
$ws.Range("J3").Copy()
$ws.Range("H3:I3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0

# --- Row 4 (Bus 2): load becomes what Bus 1 used to carry ---
$ws.Cells.Item(4, 8).Value = 1500
$ws.Cells.Item(4, 9).Value = 750

# --- Insert a new row 5 for "Bus 3" (formatting is copied down from row 4) ---
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Bus 3"
$ws.Cells.Item(5, 3).Value = "MV 1"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.9828
$ws.Cells.Item(5, 7).Value = 0.0996
$ws.Cells.Item(5, 8).Value = 900
$ws.Cells.Item(5, 9).Value = 500
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 12.66
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0

# --- Update selection to reflect new layout ---
$ws.Range("I7").Select()
